$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (updated crypto price/volume data).
$updates = @{
    "D2" = "28.202.87"
    "E2" = "  +1.49%  "
    "D3" = "1.794.23"
    "E3" = "  +2.93%  "
    "E4" = "  +0.26%  "
    "D5" = "335.61"
    "E5" = "  +0.62%  "
    "E6" = "  +0.18%  "
    "D7" = "0.4496"
    "E7" = "  +15.74%  "
    "D8" = "0.3713"
    "E8" = "  +10.08%  "
    "D9" = "45.12"
    "E9" = "  -0.53%  "
    "B10" = "Polygon"
    "C10" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D10" = "1.142"
    "E10" = "  +3.85%  "
    "B11" = "Dogecoin"
    "C11" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "D11" = "0.07565"
    "E11" = "  +5.80%  "
    "D12" = "1.003"
    "E12" = "  +0.31%  "
    "D13" = "22.29"
    "E13" = "  +2.09%  "
    "D14" = "6.288"
    "E14" = "  +3.59%  "
    "D15" = "7.445"
    "E15" = "  +7.41%  "
    "D16" = "1.793.39"
    "E16" = "  +2.95%  "
    "D17" = "0.00001088"
    "E17" = "  +3.85%  "
    "D18" = "0.06732"
    "E18" = "  +1.93%  "
    "D19" = "81.10"
    "E19" = "  +2.64%  "
    "E20" = "  +0.25%  "
    "D21" = "17.45"
    "E21" = "  +4.28%  "
    "D22" = "6.364"
    "E22" = "  +3.58%  "
    "D23" = "28.204.77"
    "E23" = "  +1.60%  "
    "D24" = "11.77"
    "E24" = "  +2.49%  "
    "D25" = "2.420"
    "E25" = "  +1.48%  "
    "D26" = "20.50"
    "E26" = "  +3.92%  "
    "E27" = "  -1.48%  "
    "D28" = "2.359"
    "E28" = "  +3.78%  "
    "D29" = "1.998.10"
    "E29" = "  +2.98%  "
    "D30" = "133.10"
    "E30" = "  +4.40%  "
    "D31" = "1.237"
    "E31" = "  -2.38%  "
    "D32" = "4.039"
    "E32" = "  -0.53%  "
    "B33" = "Filecoin"
    "C33" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D33" = "5.796"
    "E33" = "  +1.10%  "
    "B34" = "Stellar"
    "C34" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D34" = "0.09377"
    "E34" = "  +7.71%  "
    "D35" = "0.2367"
    "E35" = "  +13.82%  "
    "D36" = "12.05"
    "E36" = "  +1.00%  "
    "D37" = "0.06305"
    "E37" = "  +4.06%  "
    "D38" = "0.02327"
    "E38" = "  +3.17%  "
    "D39" = "5.207"
    "E39" = "  +2.65%  "
    "D40" = "0.6550"
    "E40" = "  +2.47%  "
    "D41" = "1.480"
    "E41" = "  -2.06%  "
    "D42" = "8.293"
    "E42" = "  +5.81%  "
    "E43" = "  +1.72%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D44" = "14.16"
    "E44" = "  +5.26%  "
    "B45" = "Frax"
    "C45" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
    "D45" = "1.002"
    "E45" = "  +0.18%  "
    "D46" = "3.831"
    "E46" = "  +0.62%  "
    "D47" = "0.6068"
    "E47" = "  +2.74%  "
    "D48" = "129.64"
    "E48" = "  +3.16%  "
    "D49" = "2.021"
    "E49" = "  +2.97%  "
    "D50" = "0.07115"
    "E50" = "  +2.80%  "
    "D51" = "1.158"
    "E51" = "  +1.57%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    # Force the value to be stored as text (matches source data which is
    # plain text, not numeric), same as the other static text columns.
    if ($addr -match "^[DE]") {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
